# Generate Report for Handoff
# Updates the localization-status workbook to reflect that the
# 5ffd0142-ccb6-4c91-9c0e-c04c855f5fef.md file is now ready for handoff
# (handoff xliff is stale vs. the latest source and needs to be regenerated).

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/990aedbfbb28a4a6e117e4037a60e4625d824cfc/e2e/5ffd0142-ccb6-4c91-9c0e-c04c855f5fef.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6becde20a3c8564b316fb15cb48c66ef0d9e7072/e2e/5ffd0142-ccb6-4c91-9c0e-c04c855f5fef.md."

# --- Overview sheet: row for 5ffd0142-ccb6-4c91-9c0e-c04c855f5fef.md (row 3) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-18 06:46:32"

# --- zh-cn sheet: row for 5ffd0142-ccb6-4c91-9c0e-c04c855f5fef.md (row 3) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("H3").Value = "2016-08-18 06:46:26"
$wsZhCn.Range("P3").Value = $errorDetail
# Column P (Error Detail) widens to fit the long message, matching column A's width.
$wsZhCn.Range("P1").ColumnWidth = $wsZhCn.Range("A1").ColumnWidth()

# --- de-de sheet: row for 5ffd0142-ccb6-4c91-9c0e-c04c855f5fef.md (row 3) ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("H3").Value = "2016-08-18 06:46:32"
$wsDeDe.Range("P3").Value = $errorDetail
# Column P (Error Detail) widens to fit the long message, matching column A's width.
$wsDeDe.Range("P1").ColumnWidth = $wsDeDe.Range("A1").ColumnWidth()
